# The speed will become faster act game continue
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 (id 5001): the key-point-score multiplier goes from 1 to 2,
# which ripples the shared formula in E20 (G20 - C20*F20) from 50 to 200.
$ws.Range("C20").Value = 2

# The active selection moves on to F24.
$ws.Range("F24").Select()

# Shrink the workbook window to match the new view size.
$excel.ActiveWindow.Width = 19350
$excel.ActiveWindow.Height = 11460
